$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.495.78"
$ws.Range("E2").Value = "  +1.45%  "
$ws.Range("D3").Value = "1.676.98"
$ws.Range("E3").Value = "  +1.67%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5311"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.69%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2697"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.47%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06413"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.45%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.80"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.77%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07809"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.53%  "
$ws.Range("D12").Value = "1.680.88"
$ws.Range("E12").Value = "  +1.91%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.510"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.03%  "
$ws.Range("E14").Value = "  +0.25%  "
$ws.Range("D15").Value = "0.0₅8354"
$ws.Range("E15").Value = "  +1.97%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.78"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.75%  "
$ws.Range("D17").Value = "26.534.52"
$ws.Range("E17").Value = "  +1.64%  "
$ws.Range("E18").Value = "  -0.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.791"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "193.13"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.33"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.88%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.321"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.70%  "
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("E24").Value = "  +5.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "138.91"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.84%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.412"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.30"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.59%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.439"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.06319"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.69%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.293"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.606"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.92%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.443"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.99%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.691"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.23%  "
$ws.Range("E34").Value = "  +2.63%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6171"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +8.71%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.426"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.36%  "
$ws.Range("E37").Value = "  +1.10%  "
$ws.Range("E38").Value = "  +1.00%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.122"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.89%  "
$ws.Range("D40").Value = "1.094.73"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8633"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.36%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.69"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.46%  "
$ws.Range("D44").Value = "1.823.60"
$ws.Range("E44").Value = "  +1.48%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "58.65"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.65%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.0₈109"
$ws.Range("E46").Value = "  +1.91%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.200"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.30%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9945"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.71%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.494"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.89%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05194"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.027"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.54%  "
